# Wijzigingen en aanvullingen tijdreeksen.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "prog001" to "data"
$ws.Name = "data"

# Update header text (drop the "(x1000 €)" suffix, values below are now
# expressed directly in euros instead of x1000 euros)
$ws.Range("B1").Value = "inkomsten"

# Narrower column now that values are whole numbers
$ws.Columns.Item(2).ColumnWidth = 9.67

# Update the income values: previously stored as x1000 euro with 3 decimals,
# now stored as plain euro whole numbers (value * 1000)
$ws.Range("B2").Value = 2000
$ws.Range("B3").Value = 2498
$ws.Range("B4").Value = 3684
$ws.Range("B5").Value = 5162
$ws.Range("B6").Value = 6612
$ws.Range("B7").Value = 8375
$ws.Range("B8").Value = 9412
$ws.Range("B9").Value = 10880
$ws.Range("B10").Value = 11799
$ws.Range("B11").Value = 13211
$ws.Range("B12").Value = 14140
$ws.Range("B13").Value = 16652
$ws.Range("B14").Value = 18709
$ws.Range("B15").Value = 19449
$ws.Range("B16").Value = 18398
$ws.Range("B17").Value = 17312
$ws.Range("B18").Value = 14186
$ws.Range("B19").Value = 7053

# Switch number format from 3-decimal to plain integer
$ws.Range("B2:B19").NumberFormat = "0"

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
